# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# plus two pairs of rows whose ranking order swapped (B/C/D/E for rows 37/38
# and 44/45). Values that look like plain numbers are written with a leading
# apostrophe so Excel stores them as literal text (matching the original
# inlineStr cell content) instead of silently parsing them into numeric
# values (e.g. "1.000" -> 1, "0.00001134" -> 1.134E-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.822.75"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.116.08"
$ws.Range("E3").Value = "  +6.38%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'332.90"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5334"
$ws.Range("E7").Value = "  +4.59%  "
$ws.Range("D8").Value = "'0.4393"
$ws.Range("E8").Value = "  +6.87%  "
$ws.Range("D9").Value = "'0.08998"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "'47.05"
$ws.Range("E10").Value = "  +10.34%  "
$ws.Range("D11").Value = "'1.182"
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("D12").Value = "'25.06"
$ws.Range("D13").Value = "2.117.41"
$ws.Range("E13").Value = "  +6.83%  "
$ws.Range("D14").Value = "'6.773"
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("D15").Value = "'7.829"
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("D16").Value = "'97.09"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "'0.00001134"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'0.06672"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "'19.17"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'6.345"
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").Value = "30.876.63"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").Value = "'12.38"
$ws.Range("E24").Value = "  +7.72%  "
$ws.Range("D25").Value = "2.367.78"
$ws.Range("E25").Value = "  +7.05%  "
$ws.Range("D26").Value = "'2.288"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").Value = "'22.83"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'2.601"
$ws.Range("E28").Value = "  +9.65%  "
$ws.Range("D29").Value = "'163.54"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'133.43"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "'1.185"
$ws.Range("E31").Value = "  +4.40%  "
$ws.Range("D32").Value = "'0.1083"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").Value = "'6.268"
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").Value = "'1.573"
$ws.Range("E35").Value = "  +19.16%  "
$ws.Range("D36").Value = "'0.02606"
$ws.Range("E36").Value = "  +5.10%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.563"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'12.91"
$ws.Range("E38").Value = "  +10.08%  "
$ws.Range("D39").Value = "'0.06776"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").Value = "'9.513"
$ws.Range("E40").Value = "  +7.02%  "
$ws.Range("D41").Value = "'0.2293"
$ws.Range("D42").Value = "'0.6875"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").Value = "'1.253"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.19"
$ws.Range("E44").Value = "  +4.40%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6468"
$ws.Range("E45").Value = "  +5.66%  "
$ws.Range("D46").Value = "'0.9996"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'2.231"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "'3.667"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "'1.279"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("D50").Value = "'83.27"
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("D51").Value = "'121.84"
$ws.Range("E51").Value = "  -1.61%  "
